$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.198.23'
$ws.Range("E2").Value = '  +3.68%  '

$ws.Range("D3").Value = '2.430.52'
$ws.Range("E3").Value = '  +0.67%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.01'
$ws.Range("E5").Value = '  +3.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.76'
$ws.Range("E6").Value = '  +5.38%  '

$ws.Range("E7").Value = '  +1.50%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.527'
$ws.Range("E9").Value = '  +7.25%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.54'
$ws.Range("E10").Value = '  +0.91%  '

$ws.Range("E11").Value = '  +0.60%  '

$ws.Range("E12").Value = '  -2.12%  '

$ws.Range("E13").Value = '  -2.05%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.04'
$ws.Range("E14").Value = '  +2.25%  '

$ws.Range("D15").Value = '2.807.68'
$ws.Range("E15").Value = '  +0.99%  '

$ws.Range("D16").Value = '2.438.57'
$ws.Range("E16").Value = '  +1.42%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.841'
$ws.Range("E17").Value = '  +1.92%  '

$ws.Range("D18").Value = '45.100.77'
$ws.Range("E18").Value = '  +3.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.28'
$ws.Range("E19").Value = '  +1.02%  '

$ws.Range("E20").Value = '  -0.98%  '

$ws.Range("D21").Value = '0.0₃0920'
$ws.Range("E21").Value = '  +2.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.86'
$ws.Range("E22").Value = '  +0.77%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '244.12'
$ws.Range("E23").Value = '  +2.53%  '

$ws.Range("E24").Value = '  +0.48%  '

$ws.Range("E25").Value = '  +1.55%  '

$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.58'
$ws.Range("E27").Value = '  +2.44%  '

$ws.Range("E28").Value = '  +1.57%  '

$ws.Range("E29").Value = '  -11.81%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.05'
$ws.Range("E30").Value = '  +1.82%  '

$ws.Range("B31").Value = 'OKB'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '49.21'
$ws.Range("E31").Value = '  +2.77%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.36'
$ws.Range("E32").Value = '  +10.35%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.125'
$ws.Range("E33").Value = '  +6.21%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.22'
$ws.Range("E34").Value = '  +1.67%  '

$ws.Range("E35").Value = '  +0.40%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0764'
$ws.Range("E36").Value = '  +1.74%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.88'
$ws.Range("E37").Value = '  -0.67%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.44'
$ws.Range("E38").Value = '  +0.80%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.86'
$ws.Range("E39").Value = '  -2.39%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '123.82'
$ws.Range("E40").Value = '  -4.99%  '

$ws.Range("E41").Value = '  -2.93%  '

$ws.Range("E42").Value = '  +1.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.57'
$ws.Range("E43").Value = '  -2.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0289'
$ws.Range("E44").Value = '  +2.01%  '

$ws.Range("D45").Value = '1.935.88'

$ws.Range("E46").Value = '  -2.88%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.93'
$ws.Range("E47").Value = '  +3.14%  '

$ws.Range("E48").Value = '  +15.88%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.19'
$ws.Range("E49").Value = '  -1.31%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '76.51'
$ws.Range("E50").Value = '  +5.79%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '53.90'
$ws.Range("E51").Value = '  +2.17%  '
